# Generate Report for Handback
#
# Row 3 (file 190181f9-32fd-4651-9c6d-1c3ba718b3a6.md) failed the handback
# transform for both target locales, so:
#   - Overview sheet: zh-cn / de-de status cells for that row flip from
#     "Ready for handoff" to "Handback transform failed"
#   - zh-cn / de-de sheets: the same row's Status cell gets the same text
#   - zh-cn / de-de sheets: the row's "Error Detail" cell (col P) gets a
#     message explaining the handback/handoff file name mismatch
#   - zh-cn / de-de sheets: the "Error Detail" column is widened to fit

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the 190181f9-... file; E = zh-cn status, F = de-de status
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn sheet: row 3 Status column (C) + Error Detail column (P)
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("P3").Value = "Handback file name: g05kfrud.35k is different with handoff file name: 190181f9-32fd-4651-9c6d-1c3ba718b3a6.8fcb92f608bc0aa2502398589274c2b3d2dc4e89.zh-cn."

# de-de sheet: row 3 Status column (C) + Error Detail column (P)
$dede.Range("C3").Value = $newStatus
$dede.Range("P3").Value = "Handback file name: g05kfrud.35k is different with handoff file name: 190181f9-32fd-4651-9c6d-1c3ba718b3a6.8fcb92f608bc0aa2502398589274c2b3d2dc4e89.de-de."

# Widen the "Error Detail" column (16th / P) on both locale sheets to fit
# the new, much longer error messages.
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
